$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set number format to text for the Price/Volume columns so numeric-looking
# strings (e.g. "7.30", "12.34") are preserved exactly as text, matching the
# original inlineStr cell contents.
$valueRange = $ws.Range("B2:E51")
$valueRange.NumberFormat = "@"

$ws.Range("D2").Value = "71.623.05"
$ws.Range("E2").Value = "  +2.23%  "
$ws.Range("D3").Value = "3.815.59"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "704.99"
$ws.Range("E5").Value = "  +6.30%  "
$ws.Range("D6").Value = "174.90"
$ws.Range("E6").Value = "  +5.06%  "
$ws.Range("D7").Value = "3.816.00"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("D11").Value = "7.30"
$ws.Range("E11").Value = "  +4.43%  "
$ws.Range("D12").Value = "0.462"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("E13").Value = "  +7.27%  "
$ws.Range("D14").Value = "36.45"
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("D15").Value = "4.456.43"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "3.815.08"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "71.492.47"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("D18").Value = "17.75"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "7.23"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "10.92"
$ws.Range("E21").Value = "  +7.38%  "
$ws.Range("D22").Value = "483.81"
$ws.Range("E22").Value = "  +2.03%  "
$ws.Range("D23").Value = "0.716"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("D24").Value = "84.60"
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("D26").Value = "12.34"
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").Value = "10.57"
$ws.Range("E27").Value = "  +2.36%  "
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("D29").Value = "3.964.79"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").Value = "3.15"
$ws.Range("E30").Value = "  +12.35%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "7.65"
$ws.Range("E32").Value = "  +4.34%  "
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("D34").Value = "0.190"
$ws.Range("E34").Value = "  +6.56%  "
$ws.Range("D35").Value = "29.64"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("E36").Value = "  +2.82%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  +1.97%  "
$ws.Range("E39").Value = "  +4.98%  "
$ws.Range("D40").Value = "6.05"
$ws.Range("E40").Value = "  +2.45%  "
$ws.Range("D41").Value = "2.28"
$ws.Range("E41").Value = "  +11.07%  "
$ws.Range("D42").Value = "0.988"
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D45").Value = "0.000317"
$ws.Range("E45").Value = "  +16.59%  "
$ws.Range("D46").Value = "165.08"
$ws.Range("E46").Value = "  +3.83%  "
$ws.Range("D47").Value = "44.99"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").Value = "48.76"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").Value = "420.81"
$ws.Range("E49").Value = "  +8.20%  "
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").Value = "0.303"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("E51").Value = "  -2.70%  "

# Restore the default (unstyled) cell style now that the text values are set,
# so no stray number-format styling is left behind on the cells.
$valueRange.Style = "Normal"
